$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.451.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.660.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.24%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  +8.09%  "

$ws.Range("E10").Value = "  +3.33%  "

$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000190"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +15.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.138.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.211.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.665.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.24%  "

$ws.Range("E21").Value = "  +4.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.04%  "

$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.18%  "

$ws.Range("E26").Value = "  +16.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("E28").Value = "  -1.86%  "

$ws.Range("E29").Value = "  +1.92%  "

$ws.Range("E30").Value = "  +6.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "545.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("E33").Value = "  +2.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.38%  "

$ws.Range("E35").Value = "  +4.27%  "

$ws.Range("E36").Value = "  +3.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("E39").Value = "  +0.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("E44").Value = "  +2.41%  "

$ws.Range("E45").Value = "  +6.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0616"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.660"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.60%  "

$ws.Range("E49").Value = "  +4.33%  "

$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.83%  "
